# Fill in "Day 5" (column F) of the third week's sleep-diary table
# (rows 44-57), which was previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Time-of-day entries (rows 44-47): wake time, get-up time, bed time,
# --- lights-off time. These need the "h:mm" time number format, same as
# --- the other day columns in this table.
$ws.Range("F44:F47").NumberFormat = "h:mm"

$ws.Range("F44").Value = 0.32083333333333336   # 7:42
$ws.Range("F45").Value = 0.34027777777777779   # 8:10
$ws.Range("F46").Value = 0.96180555555555558   # 23:05
$ws.Range("F47").Value = 0.97222222222222221   # 23:20

# --- Plain numeric entries ---
$ws.Range("F48").Value = 20
$ws.Range("F49").Value = 2
$ws.Range("F50").Value = 5
$ws.Range("F51").Value = 480

# --- Text entries (matches the "无" used in column E for the same rows) ---
$ws.Range("F52").Value = "无"

$ws.Range("F53").Value = 35
$ws.Range("F54").Value = 3
$ws.Range("F55").Value = 1
$ws.Range("F56").Value = 2

$ws.Range("F57").Value = "无"

# --- Update the sheet's active cell / selection to match the final state ---
$null = $ws.Range("F57").Select()
